$d = $word.ActiveDocument

# The title paragraph currently starts with a single run containing
# "Level 3 " (trailing space) immediately followed by "Teleportation ".
# The target splits that run right after "Level 3" and drops a
# "_GoBack" bookmark into the gap, leaving the trailing space as its own
# run:
#   <w:r><w:t>Level 3</w:t></w:r>
#   <w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>
#   <w:r><w:t xml:space="preserve"> </w:t></w:r>
#   <w:r><w:t xml:space="preserve">Teleportation </w:t></w:r>   (untouched)
#
# Word keeps bookmark names unique, so (re)adding a bookmark named
# "_GoBack" here also removes the pre-existing "_GoBack" bookmark that
# currently lives a few paragraphs further down - which is exactly the
# second half of the requested change.

# 1) Split "Level 3 " into "Level 3" + " " by dropping the bookmark
#    right after the word "3" (position 7, measured from the start of
#    the document/paragraph).
$splitPoint = $d.Range(7, 7)
$d.Bookmarks.Add("_GoBack", $splitPoint)

# 2) The split above leaves the xml:space="preserve" flag on the new
#    "Level 3" run (it no longer needs it, since it has no leading or
#    trailing whitespace). Touch the run's text (swap it out for a
#    same-length placeholder and back) so the run gets re-serialized and
#    the now-unneeded xml:space="preserve" is dropped, without touching
#    any text outside the bookmark (the bookmark stops the
#    re-serialization from spilling into the following
#    "Teleportation "/"Blueprint:" runs).
$levelRange = $d.Range(0, 7)
$levelRange.Text = "Level X"
$levelRange2 = $d.Range(0, 7)
$levelRange2.Text = "Level 3"
